$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the _GoBack bookmark from its current location (start of
#    the "AnyOject / 代表任何的类对象" heading paragraph).
# ------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# ------------------------------------------------------------------
# 2. Locate the paragraph that ends with "...表示任意的类型，包括整型，浮点型"
#    (the "Any" heading) - the four new paragraphs are inserted right
#    after it (and before the pre-existing blank paragraph that
#    follows it).
# ------------------------------------------------------------------
$anchorRange = $d.Content
$null = $anchorRange.Find.Execute("表示任意的类型，包括整型，浮点型", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$containingPara = $anchorRange.Paragraphs(1)
$insertAt = $d.Range($containingPara.Range.End, $containingPara.Range.End)

# ------------------------------------------------------------------
# 3. Build the OOXML for the four new paragraphs:
#      6. 扩展 extension
#      ？嵌套类型没看
#      (empty paragraph)
#      7.协议   <- gets the _GoBack bookmark back, at its end
#    A final throw-away paragraph is appended so that Word's "last
#    paragraph merges into the following paragraph" splice behaviour
#    consumes a disposable paragraph instead of corrupting the
#    pre-existing blank paragraph that must stay untouched; the
#    left-over marker text is deleted again right after insertion.
# ------------------------------------------------------------------
$marker = "ZZ_TMP_MARKER_ZZ"

$fragment = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:body>' +
            '<w:p>' +
              '<w:pPr><w:pStyle w:val="1"/><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr>' +
              '<w:r><w:t>6.</w:t></w:r>' +
              '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>扩展</w:t></w:r>' +
              '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> extension</w:t></w:r>' +
            '</w:p>' +
            '<w:p>' +
              '<w:pPr><w:pStyle w:val="2"/><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr>' +
              '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>？</w:t></w:r>' +
              '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>嵌套类型没看</w:t></w:r>' +
            '</w:p>' +
            '<w:p>' +
              '<w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr>' +
            '</w:p>' +
            '<w:p>' +
              '<w:pPr><w:pStyle w:val="1"/><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr>' +
              '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>7.</w:t></w:r>' +
              '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>协议</w:t></w:r>' +
            '</w:p>' +
            '<w:p><w:r><w:t>' + $marker + '</w:t></w:r></w:p>' +
          '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'

$insertAt.InsertXML($fragment)

# ------------------------------------------------------------------
# 4. Strip the throw-away marker text back out again (it landed,
#    merged, inside the pre-existing blank paragraph which must stay
#    exactly as it was).
# ------------------------------------------------------------------
$cleanup = $d.Content
$null = $cleanup.Find.Execute($marker, $false, $false, $false, $false, $false, $true, 1, $false, "", 2)

# ------------------------------------------------------------------
# 5. Re-create the _GoBack bookmark at the end of the new "7.协议"
#    paragraph (matches the moved bookmarkStart/bookmarkEnd in the
#    diff).
# ------------------------------------------------------------------
$target = $d.Content
$null = $target.Find.Execute("7.协议", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bmPoint = $d.Range($target.End, $target.End)
$null = $d.Bookmarks.Add("_GoBack", $bmPoint)

Write-Output "ok"
